$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "28.290.17"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.801.24"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'326.54"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.4423"
$ws.Range("E7").Value = "  +11.62%  "
$ws.Range("D8").Value = "'0.3720"
$ws.Range("E8").Value = "  +8.72%  "
$ws.Range("D9").Value = "'44.58"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "'1.150"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.07511"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("D12").Value = "'22.55"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'1.002"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.658"
$ws.Range("E14").Value = "  +6.10%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'6.297"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "1.799.34"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "'0.00001091"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "'0.06760"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'80.99"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'17.47"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").Value = "'6.327"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "28.288.05"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "'11.79"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'2.419"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'20.40"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'152.94"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").Value = "'2.367"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "1.997.31"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").Value = "'132.67"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'1.231"
$ws.Range("E31").Value = "  -4.53%  "
$ws.Range("D32").Value = "'4.025"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").Value = "'5.824"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "'0.09347"
$ws.Range("E34").Value = "  +6.48%  "
$ws.Range("D35").Value = "'0.2313"
$ws.Range("E35").Value = "  +8.34%  "
$ws.Range("D36").Value = "'12.11"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "'0.06327"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D39").Value = "'5.156"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'0.6565"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'1.477"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("D42").Value = "'8.192"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "'1.202"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'14.00"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "'0.6060"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'3.789"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "'129.12"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "'2.035"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").Value = "'0.07123"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "'1.156"
$ws.Range("E51").Value = "  -0.75%  "
